# mgt_fee/Quarterly_mgt_fee/2022Q1Mgmt.xlsx - "add mgt fee ppt"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 10 data corrections:
#  - H10 no longer carries a management-fee number (cleared)
#  - J10 becomes a literal 0 instead of the blank shared string
$ws.Range("H10").ClearContents()
$ws.Range("J10").Value = 0

# Selection moves from C9 down onto the whole of row 10
$ws.Range("A10:J10").Select()
